$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($cellRef).Style = 'Normal'
}

Set-TextValue 'D2' '42.651.96'
Set-TextValue 'E2' '  -0.74%  '
Set-TextValue 'D3' '2.527.96'
Set-TextValue 'E3' '  -1.88%  '
Set-TextValue 'D4' '0.999'
Set-TextValue 'E4' '  -0.04%  '
Set-TextValue 'D5' '308.51'
Set-TextValue 'E5' '  -2.10%  '
Set-TextValue 'D6' '100.57'
Set-TextValue 'E6' '  +0.61%  '
Set-TextValue 'E7' '  -1.32%  '
Set-TextValue 'E8' '  +0.06%  '
Set-TextValue 'E9' '  -2.79%  '
Set-TextValue 'D10' '35.71'
Set-TextValue 'E10' '  -1.49%  '
Set-TextValue 'D11' '0.0804'
Set-TextValue 'E11' '  -1.22%  '
Set-TextValue 'D12' '7.36'
Set-TextValue 'E12' '  -2.46%  '
Set-TextValue 'E13' '  -0.02%  '
Set-TextValue 'D14' '2.917.70'
Set-TextValue 'E14' '  -1.85%  '
Set-TextValue 'B15' 'WrappedEther'
Set-TextValue 'C15' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D15' '2.542.86'
Set-TextValue 'E15' '  -3.81%  '
Set-TextValue 'B16' 'Chainlink'
Set-TextValue 'C16' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D16' '15.28'
Set-TextValue 'E16' '  -2.79%  '
Set-TextValue 'E17' '  -4.26%  '
Set-TextValue 'D18' '42.636.73'
Set-TextValue 'E18' '  -0.89%  '
Set-TextValue 'D19' '6.70'
Set-TextValue 'E19' '  -2.41%  '
Set-TextValue 'D20' '0.0₃0949'
Set-TextValue 'E20' '  -2.10%  '
Set-TextValue 'D21' '12.21'
Set-TextValue 'E21' '  -4.09%  '
Set-TextValue 'D22' '69.55'
Set-TextValue 'E22' '  +0.11%  '
Set-TextValue 'D23' '242.74'
Set-TextValue 'E23' '  -2.95%  '
Set-TextValue 'E24' '  -3.46%  '
Set-TextValue 'E25' '  -2.90%  '
Set-TextValue 'E26' '  -0.02%  '
Set-TextValue 'D27' '25.42'
Set-TextValue 'E27' '  -6.28%  '
Set-TextValue 'E28' '  -2.78%  '
Set-TextValue 'D29' '10.12'
Set-TextValue 'E29' '  -2.00%  '
Set-TextValue 'D30' '38.38'
Set-TextValue 'E30' '  -5.51%  '
Set-TextValue 'D31' '157.37'
Set-TextValue 'E31' '  -0.01%  '
Set-TextValue 'E32' '  -1.54%  '
Set-TextValue 'E33' '  +10.00%  '
Set-TextValue 'E34' '  -1.57%  '
Set-TextValue 'E35' '  -2.82%  '
Set-TextValue 'B36' 'Celestia'
Set-TextValue 'C36' 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue 'D36' '18.00'
Set-TextValue 'E36' '  -4.04%  '
Set-TextValue 'B37' 'LidoDAOToken'
Set-TextValue 'C37' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D37' '3.13'
Set-TextValue 'E37' '  -8.82%  '
Set-TextValue 'D38' '1.97'
Set-TextValue 'E38' '  -7.40%  '
Set-TextValue 'E39' '  -1.67%  '
Set-TextValue 'E40' '  -1.02%  '
Set-TextValue 'E41' '  +2.28%  '
Set-TextValue 'D42' '22.21'
Set-TextValue 'E42' '  -6.15%  '
Set-TextValue 'E43' '  +0.02%  '
Set-TextValue 'E44' '  -1.43%  '
Set-TextValue 'D45' '3.28'
Set-TextValue 'E45' '  +0.62%  '
Set-TextValue 'D46' '2.003.64'
Set-TextValue 'E46' '  +0.01%  '
Set-TextValue 'E47' '  -0.73%  '
Set-TextValue 'D48' '2.771.72'
Set-TextValue 'E48' '  -1.84%  '
Set-TextValue 'D49' '0.189'
Set-TextValue 'E49' '  -4.41%  '
Set-TextValue 'D50' '79.13'
Set-TextValue 'E50' '  -3.47%  '
Set-TextValue 'D51' '71.90'
Set-TextValue 'E51' '  -4.39%  '
